# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the c84cf48f-f4f9-4261-9447-c117d69a89ae.md row (row 5) to reflect a new
# handoff report generation pass.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G5").Value = "2016-11-09 00:07:48"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H5").Value = "2016-11-09 00:07:35"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H5").Value = "2016-11-09 00:07:48"
